$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update value "TDI" -> "InADvance" for the ITSM/Copec row
$ws.Range("C5").Value = "InADvance"

# Rename header "Categoria" -> "Partner"
$ws.Range("C1").Value = "Partner"

# Update the active selection to C1 (single cell)
$ws.Range("C1").Select()
